# Add three new line items to the "Instal days by Model" table, keeping
# the existing alphabetical sort of the Item column, then leave that sheet
# as the active one with C40 selected (matching the authored workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instal days by Model")

# --- Insert "Production Support Day" just above "R225" (row 23) ---
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "Production Support Day"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 0

# --- Insert "Train Day - Eng" and "Train Day - Tech" just above "UC15" ---
# "UC15" was row 37 before the first insert above pushed everything below
# row 23 down by one, so it is now row 38.
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "Train Day - Eng"
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 1

$ws.Rows.Item(39).Insert()
$ws.Range("A39").Value = "Train Day - Tech"
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 0

# --- Grow Table1 (and its autofilter) to cover the three new rows ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:C40"))

# Keep the workbook-level hidden _FilterDatabase name (driven by the
# table's AutoFilter) in sync with the table's new extent.
$filterName = $wb.Names.Item("Instal days by Model!_FilterDatabase")
$filterName.RefersTo = "='Instal days by Model'!`$A`$1:`$C`$40"

# --- Match the saved selection/active-sheet state from the edit ---
$ws.Activate()
$ws.Range("C40").Select()
